# Generate Report for Handoff
# Adds two newly-handed-off source files
#   25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md
#   65c0ee2c-438f-4b42-b24f-0f480980a8bf.md
# to the localization status report (Overview + zh-cn + de-de sheets),
# pushing the existing d4502247-... row down to the bottom of each table.

function Set-CellText {
    param($range, [string]$text)
    # Excel auto-coerces the bare words True/False into boolean cells.
    # Force them (and anything else that looks ambiguous) to stay text by
    # using the classic leading-apostrophe text prefix; the apostrophe
    # itself never lands in the stored value.
    if ($text -eq "True" -or $text -eq "False") {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 3: new entry 25aebf7b...
Set-CellText $ov.Range("A3") "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md"
Set-CellText $ov.Range("B3") 'e2e\25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md'
Set-CellText $ov.Range("C3") ".md"
Set-CellText $ov.Range("D3") ""
Set-CellText $ov.Range("E3") "Ready for handoff"
Set-CellText $ov.Range("F3") "Ready for handoff"
Set-CellText $ov.Range("G3") "2016-08-26 04:39:23"

# Row 4: new entry 65c0ee2c...
Set-CellText $ov.Range("A4") "65c0ee2c-438f-4b42-b24f-0f480980a8bf.md"
Set-CellText $ov.Range("B4") 'e2e\65c0ee2c-438f-4b42-b24f-0f480980a8bf.md'
Set-CellText $ov.Range("C4") ".md"
Set-CellText $ov.Range("D4") ""
Set-CellText $ov.Range("E4") "Ready for handoff"
Set-CellText $ov.Range("F4") "Ready for handoff"
Set-CellText $ov.Range("G4") "2016-08-26 04:39:23"

# Row 5: existing d4502247... entry, now pushed to the bottom
Set-CellText $ov.Range("A5") "d4502247-d3c0-4021-bb47-4c29efc6528d.md"
Set-CellText $ov.Range("B5") 'e2e\d4502247-d3c0-4021-bb47-4c29efc6528d.md'
Set-CellText $ov.Range("C5") ".md"
Set-CellText $ov.Range("D5") ""
Set-CellText $ov.Range("E5") "Ready for handoff"
Set-CellText $ov.Range("F5") "Ready for handoff"
Set-CellText $ov.Range("G5") "2016-08-26 04:38:12"

$ov.Range("B3").Style = "HyperLink"
$ov.Range("B4").Style = "HyperLink"
$ov.Range("B5").Style = "HyperLink"

$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md", "", "", 'e2e\25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md') | Out-Null
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/65c0ee2c-438f-4b42-b24f-0f480980a8bf.md", "", "", 'e2e\65c0ee2c-438f-4b42-b24f-0f480980a8bf.md') | Out-Null
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd09a6b63850306670958e7b4e0c8d8e9aafa722/e2e/d4502247-d3c0-4021-bb47-4c29efc6528d.md", "", "", 'e2e\d4502247-d3c0-4021-bb47-4c29efc6528d.md') | Out-Null

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G5"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 3: new entry 25aebf7b...
Set-CellText $zh.Range("A3") "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md"
Set-CellText $zh.Range("B3") ".md"
Set-CellText $zh.Range("C3") "Ready for handoff"
Set-CellText $zh.Range("D3") "e2e"
Set-CellText $zh.Range("E3") "ht"
Set-CellText $zh.Range("F3") "False"
Set-CellText $zh.Range("G3") "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.8af9a6956460890549a3c0ccafd5892e7936d779.zh-cn.xlf"
Set-CellText $zh.Range("H3") "2016-08-26 04:39:18"
Set-CellText $zh.Range("I3") ""
Set-CellText $zh.Range("J3") ""
Set-CellText $zh.Range("K3") "0001-01-01 00:00:00"
Set-CellText $zh.Range("L3") ""
Set-CellText $zh.Range("M3") "True"
Set-CellText $zh.Range("N3") ""
Set-CellText $zh.Range("O3") "False"
Set-CellText $zh.Range("P3") ""

# Row 4: new entry 65c0ee2c...
Set-CellText $zh.Range("A4") "65c0ee2c-438f-4b42-b24f-0f480980a8bf.md"
Set-CellText $zh.Range("B4") ".md"
Set-CellText $zh.Range("C4") "Ready for handoff"
Set-CellText $zh.Range("D4") "e2e"
Set-CellText $zh.Range("E4") "ht"
Set-CellText $zh.Range("F4") "False"
Set-CellText $zh.Range("G4") "65c0ee2c-438f-4b42-b24f-0f480980a8bf.482da8da870f00ed602e2a5e5d4444fe0901f789.zh-cn.xlf"
Set-CellText $zh.Range("H4") "2016-08-26 04:39:18"
Set-CellText $zh.Range("I4") ""
Set-CellText $zh.Range("J4") ""
Set-CellText $zh.Range("K4") "0001-01-01 00:00:00"
Set-CellText $zh.Range("L4") ""
Set-CellText $zh.Range("M4") "True"
Set-CellText $zh.Range("N4") ""
Set-CellText $zh.Range("O4") "False"
Set-CellText $zh.Range("P4") ""

# Row 5: existing d4502247... entry, now pushed to the bottom
Set-CellText $zh.Range("A5") "d4502247-d3c0-4021-bb47-4c29efc6528d.md"
Set-CellText $zh.Range("B5") ".md"
Set-CellText $zh.Range("C5") "Ready for handoff"
Set-CellText $zh.Range("D5") "e2e"
Set-CellText $zh.Range("E5") "ht"
Set-CellText $zh.Range("F5") "False"
Set-CellText $zh.Range("G5") "d4502247-d3c0-4021-bb47-4c29efc6528d.7629d3f4c35f598ade6d53ad933e5aa550516ef1.zh-cn.xlf"
Set-CellText $zh.Range("H5") "2016-08-26 04:38:07"
Set-CellText $zh.Range("I5") ""
Set-CellText $zh.Range("J5") ""
Set-CellText $zh.Range("K5") "0001-01-01 00:00:00"
Set-CellText $zh.Range("L5") ""
Set-CellText $zh.Range("M5") "True"
Set-CellText $zh.Range("N5") ""
Set-CellText $zh.Range("O5") "False"
Set-CellText $zh.Range("P5") ""

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md", "", "", "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/65c0ee2c-438f-4b42-b24f-0f480980a8bf.md", "", "", "65c0ee2c-438f-4b42-b24f-0f480980a8bf.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd09a6b63850306670958e7b4e0c8d8e9aafa722/e2e/d4502247-d3c0-4021-bb47-4c29efc6528d.md", "", "", "d4502247-d3c0-4021-bb47-4c29efc6528d.md") | Out-Null

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P5"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 3: new entry 25aebf7b...
Set-CellText $de.Range("A3") "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md"
Set-CellText $de.Range("B3") ".md"
Set-CellText $de.Range("C3") "Ready for handoff"
Set-CellText $de.Range("D3") "e2e"
Set-CellText $de.Range("E3") "ht"
Set-CellText $de.Range("F3") "False"
Set-CellText $de.Range("G3") "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.8af9a6956460890549a3c0ccafd5892e7936d779.de-de.xlf"
Set-CellText $de.Range("H3") "2016-08-26 04:39:23"
Set-CellText $de.Range("I3") ""
Set-CellText $de.Range("J3") ""
Set-CellText $de.Range("K3") "0001-01-01 00:00:00"
Set-CellText $de.Range("L3") ""
Set-CellText $de.Range("M3") "True"
Set-CellText $de.Range("N3") ""
Set-CellText $de.Range("O3") "False"
Set-CellText $de.Range("P3") ""

# Row 4: new entry 65c0ee2c...
Set-CellText $de.Range("A4") "65c0ee2c-438f-4b42-b24f-0f480980a8bf.md"
Set-CellText $de.Range("B4") ".md"
Set-CellText $de.Range("C4") "Ready for handoff"
Set-CellText $de.Range("D4") "e2e"
Set-CellText $de.Range("E4") "ht"
Set-CellText $de.Range("F4") "False"
Set-CellText $de.Range("G4") "65c0ee2c-438f-4b42-b24f-0f480980a8bf.482da8da870f00ed602e2a5e5d4444fe0901f789.de-de.xlf"
Set-CellText $de.Range("H4") "2016-08-26 04:39:23"
Set-CellText $de.Range("I4") ""
Set-CellText $de.Range("J4") ""
Set-CellText $de.Range("K4") "0001-01-01 00:00:00"
Set-CellText $de.Range("L4") ""
Set-CellText $de.Range("M4") "True"
Set-CellText $de.Range("N4") ""
Set-CellText $de.Range("O4") "False"
Set-CellText $de.Range("P4") ""

# Row 5: existing d4502247... entry, now pushed to the bottom
Set-CellText $de.Range("A5") "d4502247-d3c0-4021-bb47-4c29efc6528d.md"
Set-CellText $de.Range("B5") ".md"
Set-CellText $de.Range("C5") "Ready for handoff"
Set-CellText $de.Range("D5") "e2e"
Set-CellText $de.Range("E5") "ht"
Set-CellText $de.Range("F5") "False"
Set-CellText $de.Range("G5") "d4502247-d3c0-4021-bb47-4c29efc6528d.7629d3f4c35f598ade6d53ad933e5aa550516ef1.de-de.xlf"
Set-CellText $de.Range("H5") "2016-08-26 04:38:12"
Set-CellText $de.Range("I5") ""
Set-CellText $de.Range("J5") ""
Set-CellText $de.Range("K5") "0001-01-01 00:00:00"
Set-CellText $de.Range("L5") ""
Set-CellText $de.Range("M5") "True"
Set-CellText $de.Range("N5") ""
Set-CellText $de.Range("O5") "False"
Set-CellText $de.Range("P5") ""

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md", "", "", "25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/65c0ee2c-438f-4b42-b24f-0f480980a8bf.md", "", "", "65c0ee2c-438f-4b42-b24f-0f480980a8bf.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd09a6b63850306670958e7b4e0c8d8e9aafa722/e2e/d4502247-d3c0-4021-bb47-4c29efc6528d.md", "", "", "d4502247-d3c0-4021-bb47-4c29efc6528d.md") | Out-Null

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P5"))

Write-Output "Report regenerated: added 25aebf7b-15ad-45f1-8c2f-938b2fcfad08.md and 65c0ee2c-438f-4b42-b24f-0f480980a8bf.md; d4502247-... moved to the last row on Overview, zh-cn and de-de sheets."
